$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-set A1 (same text) - triggers new shared-string entries per the diff
$ws.Range("A1").Value = "HK_G_acc_LT"

# Update the accuracy values in A2:A49
$ws.Range("A2").Value = 75.246132208157519
$ws.Range("A3").Value = 74.542897327707465
$ws.Range("A4").Value = 77.215189873417728
$ws.Range("A5").Value = 67.510548523206751
$ws.Range("A6").Value = 70.745428973277072
$ws.Range("A7").Value = 71.589310829817151
$ws.Range("A8").Value = 82.278481012658233
$ws.Range("A9").Value = 82.278481012658233
$ws.Range("A10").Value = 81.57524613220815
$ws.Range("A11").Value = 82.137834036568208
$ws.Range("A12").Value = 65.682137834036567
$ws.Range("A13").Value = 70.182841068917028
$ws.Range("A14").Value = 76.933895921237692
$ws.Range("A15").Value = 76.793248945147667
$ws.Range("A16").Value = 77.637130801687761
$ws.Range("A17").Value = 60.900140646976084
$ws.Range("A18").Value = 65.119549929676509
$ws.Range("A19").Value = 67.791842475386773
$ws.Range("A20").Value = 79.043600562587898
$ws.Range("A21").Value = 80.450070323488049
$ws.Range("A22").Value = 79.324894514767934
$ws.Range("A23").Value = 64.697609001406477
$ws.Range("A24").Value = 63.009845288326304
$ws.Range("A25").Value = 61.74402250351617
$ws.Range("A26").Value = 75.386779184247544
$ws.Range("A27").Value = 72.292545710267234
$ws.Range("A28").Value = 73.55836849507736
$ws.Range("A29").Value = 72.855133614627292
$ws.Range("A30").Value = 69.338959212376935
$ws.Range("A31").Value = 72.433192686357245
$ws.Range("A32").Value = 75.246132208157519
$ws.Range("A33").Value = 75.949367088607602
$ws.Range("A34").Value = 75.949367088607602
$ws.Range("A35").Value = 69.057665260196913
$ws.Range("A36").Value = 70.604781997187061
$ws.Range("A37").Value = 59.634317862165965
$ws.Range("A38").Value = 71.167369901547119
$ws.Range("A39").Value = 67.510548523206751
$ws.Range("A40").Value = 67.932489451476798
$ws.Range("A41").Value = 73.277074542897324
$ws.Range("A42").Value = 77.355836849507725
$ws.Range("A43").Value = 75.808720112517577
$ws.Range("A44").Value = 75.386779184247544
$ws.Range("A45").Value = 77.49648382559775
$ws.Range("A46").Value = 76.65260196905767
$ws.Range("A47").Value = 67.36990154711674
$ws.Range("A48").Value = 64.838255977496488
$ws.Range("A49").Value = 72.714486638537267
